# Generate Report for Handoff
#
# Refresh the "Latest Handoff Datetime" (column D) for every file row that
# was just re-handed-off. Row 5 (currently "In Translation") and row 11
# (the .localization-config / "Ignored" row) are left untouched since they
# were not part of this handoff run.

$wb = $excel.ActiveWorkbook

$rows = @(4, 6, 7, 8, 9, 10)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 4).Value = "2016-03-04 11:36:47"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 4).Value = "2016-03-04 11:37:01"
}
